$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Squad Total" row (row 35) contents, keep formatting
$ws.Range("A35:U35").ClearContents()

# Update selection to match the new active cell
$ws.Range("AE14").Select()
